$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Row 7: add a Person (F7) ---
$ws.Range("F7").Value = "Moors"

# --- Row 8: Klassendiagramm erstellen now has a real date/time/progress/person ---
$ws.Range("C8").Value = 42872
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("D8").Value = "11.10 Uhr - "
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial($xlPasteFormats)
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = "Moors"

# --- Row 11: Datenbankstruktur bestimmen finished ---
$ws.Range("B11").Value = "Fertig"
$ws.Range("C11").Value = 42865
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Gruppe"

# --- Row 12: Datenbank erstellen in progress ---
$ws.Range("B12").Value = "In Bearbeitung"
$ws.Range("C12").Value = 42872
$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial($xlPasteFormats)
$ws.Range("E12").Value = 0.2
$ws.Range("F12").Value = "Tahta"

# --- Row 14: Pflichtenheft bearbeiten progress update, person simplified ---
$ws.Range("E14").Value = 0.6
$ws.Range("F14").Value = "Tahta"

# --- Row 17 (new): Anwendungsfalldiagramm ---
$ws.Range("A17").Value = "Anwendungsfalldiagramm"
$ws.Range("B17").Value = "Fertig"
$ws.Range("C17").Value = 42872
$ws.Range("C2").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("D17").Value = "9.35 Uhr - 11.00 Uhr"
$ws.Range("E17").Value = 1
$ws.Range("E2").Copy()
$ws.Range("E17").PasteSpecial($xlPasteFormats)
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "Moors"

# --- Row 18 (new): Layout Design implementieren ---
$ws.Range("A18").Value = "Layout Design implementieren"
$ws.Range("B18").Value = "in Bearbeitung"
$ws.Range("C18").Value = 42872
$ws.Range("C2").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("D18").Value = "9.35 Uhr - 11.00 Uhr"
$ws.Range("D18").NumberFormat = "mmm-yy"
$ws.Range("E18").Value = 0.2
$ws.Range("E2").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("E18").Value = 0.2
$ws.Range("F18").Value = "Derksen"

# --- Row 19 (new): Backend implementieren ---
$ws.Range("A19").Value = "Backend implementieren"
$ws.Range("B19").Value = "In Bearbeitung"
$ws.Range("C19").Value = 42872
$ws.Range("C2").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)
$ws.Range("D19").Value = "9.35 Uhr - 11.00 Uhr"
$ws.Range("E19").Value = 0.2
$ws.Range("E2").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)
$ws.Range("E19").Value = 0.2
$ws.Range("F19").Value = "Horstmann"

# --- Update the view's active selection to match the final state ---
$ws.Range("E12").Select()
